# corona_cases_ksa.xlsx update
# - add a new shared string "Ras Tanura"
# - fix Tabouk -> Tabuk typo on row 109
# - add one extra confirmed case to row 183 (Jeddah / Makkah, Makkah)
# - append a new day (2020-03-30, serial 43919) of per-city rows (191-202)
#   and a trailing blank row (204)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) B109: "Tabouk" -> "Tabuk" (align the city name with the map_name H109)
# ---------------------------------------------------------------------
$ws.Range("B109").Value = "Tabuk"

# ---------------------------------------------------------------------
# 2) Row 183: new_cases 12 -> 13 (cumulative column D recalculates via
#    the existing shared formula already occupying D137:D190)
# ---------------------------------------------------------------------
$ws.Range("E183").Value = 13

# ---------------------------------------------------------------------
# 3) Row 190 gets overwritten with a brand-new record, and rows
#    191-202 are appended as new rows for 2020-03-30 (serial 43919).
# ---------------------------------------------------------------------

# Make sure every new date cell (column A) carries the same date
# formatting (style) as the rest of the column. Copy format only so we
# reuse the existing style index instead of minting a new numFmt.
$ws.Range("A189").Copy() | Out-Null
$ws.Range("A190:A202").PasteSpecial(-4122) | Out-Null
$ws.Range("A204").PasteSpecial(-4122) | Out-Null

$rows = @(
    @{ R=190; A=43919; B="Riyadh";          C="Riyadh";            E=27; H="Ar Riyad" },
    @{ R=191; A=43919; B="Dammam";          C="Eastern province";  E=14; H="Ash Sharqiyah" },
    @{ R=192; A=43919; B="Medinah";         C="Medinah";           E=23; H="Al Madinah" },
    @{ R=193; A=43919; B="Jeddah";          C="Makkah";            E=12; H="Makkah" },
    @{ R=194; A=43919; B="Mecca";           C="Makkah";            E=7;  H="Makkah" },
    @{ R=195; A=43919; B="Kobar";           C="Eastern province";  E=4;  H="Ash Sharqiyah" },
    @{ R=196; A=43919; B="Dhahran";         C="Eastern province";  E=2;  H="Ash Sharqiyah" },
    @{ R=197; A=43919; B="Qatif";           C="Eastern province";  E=2;  H="Ash Sharqiyah" },
    @{ R=198; A=43919; B="Ras Tanura";      C="Eastern province";  E=1;  H="Ash Sharqiyah" },
    @{ R=199; A=43919; B="Altaif";          C="Makkah";            E=1;  H="Makkah" },
    @{ R=200; A=43919; B="Ahsaa";           C="Eastern province";  E=1;  H="Ash Sharqiyah" },
    @{ R=201; A=43919; B="Khamis Mushait";  C="Asir";              E=1;  H='`Asir' },
    @{ R=202; A=43919; B="Tabuk";           C="Tabuk";             E=1;  H="Tabouk" }
)

foreach ($row in $rows) {
    $r = $row.R
    $ws.Range("A$r").Value = $row.A
    $ws.Range("B$r").Value = $row.B
    $ws.Range("C$r").Value = $row.C
    $prev = $r - 1
    $ws.Range("D$r").Formula = "=D$prev+E$r"
    $ws.Range("E$r").Value = $row.E
    $ws.Range("H$r").Value = $row.H
}

# Row 199 in the source workbook carries a slightly different cell
# style (explicit black font) than the rest of the date column -
# replicate it by pulling that format from the other place in the
# sheet that already uses it.
$ws.Range("M67").Copy() | Out-Null
$ws.Range("A199").PasteSpecial(-4122) | Out-Null
$ws.Range("A199").Value = 43919

# ---------------------------------------------------------------------
# 4) Trailing empty row 204 (just keeps the date-formatted style, no
#    value) - already formatted above via the PasteSpecial call.
# ---------------------------------------------------------------------

# ---------------------------------------------------------------------
# 5) Update the view so the active window matches the end of the
#    newly-appended data.
# ---------------------------------------------------------------------
$ws.Range("A203").Select() | Out-Null
